$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# New RMA record group "ZYC7" replaces the previous "NDPI" group values
# in the grid. Order mirrors the shared-string layout of the source
# workbook (RMA numbers first, then the id/name pairs per line).
$ws.Range("E2").Value = "RMA-ZYC7-001"
$ws.Range("E3").Value = "RMA-ZYC7-002"
$ws.Range("E4").Value = "RMA-ZYC7-003"

$ws.Range("J2").Value = "a7s5f000000xMriAAE"
$ws.Range("F2").Value = "RMA-ZYC7-1-1"

$ws.Range("J3").Value = "a7s5f000000xMrjAAE"
$ws.Range("F3").Value = "RMA-ZYC7-1-2"

$ws.Range("J4").Value = "a7s5f000000xMrkAAE"
$ws.Range("F4").Value = "RMA-ZYC7-1-3"
